# Add motorbikes freight data: AVLRaPTC!C7 picks up the HDV freight multiplier
# (same value as AVLRaPTC!C3) instead of a hard-coded 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AVLRaPTC")

$ws.Range("C7").Formula = "=C3"

# Reflect the author's on-screen state: AVLRaPTC tab active with C8 selected
# (the About sheet had been the active/selected tab before this edit).
$ws.Activate()
[void]$ws.Range("C8").Select()
